# Update Efna3-Epha2 NATMI TPM results with the newly-computed TPM-based
# specificity / weight figures (per commit "update scripts wuth new tpm").
# Target cluster (column D) text per row is unchanged; only the derived
# numeric columns F,G,H,M,N,O,P,Q,R,S,T differ row by row.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 (Target cluster: ECs)
$ws.Range("F2").Value = 0.3333333333333333
$ws.Range("G2").Value = 0.02551366666666667
$ws.Range("H2").Value = 0.076541
$ws.Range("M2").Value = 12.997753
$ws.Range("N2").Value = 38.993259
$ws.Range("O2").Value = 0.4740421406233454
$ws.Range("P2").Value = 0.5546503645614554
$ws.Range("Q2").Value = 0.3316203374576667
$ws.Range("R2").Value = 2.984583037119
$ws.Range("S2").Value = 0.4740421406233454
$ws.Range("T2").Value = 0.5546503645614554

# Row 3 (Target cluster: FAPs)
$ws.Range("F3").Value = 0.3333333333333333
$ws.Range("G3").Value = 0.02551366666666667
$ws.Range("H3").Value = 0.076541
$ws.Range("O3").Value = 0.02725306609819269
$ws.Range("P3").Value = 0.03188729809316786
$ws.Range("Q3").Value = 0.01906512143488889
$ws.Range("R3").Value = 0.171586092914
$ws.Range("S3").Value = 0.02725306609819269
$ws.Range("T3").Value = 0.03188729809316786

# Row 4 (Target cluster: Inflammatory-Mac)
$ws.Range("F4").Value = 0.3333333333333333
$ws.Range("G4").Value = 0.02551366666666667
$ws.Range("H4").Value = 0.076541
$ws.Range("M4").Value = 1.182591666666666
$ws.Range("N4").Value = 3.547775
$ws.Range("O4").Value = 0.04313039993528083
$ws.Range("P4").Value = 0.05046448405689858
$ws.Range("Q4").Value = 0.03017224958611111
$ws.Range("R4").Value = 0.271550246275
$ws.Range("S4").Value = 0.04313039993528083
$ws.Range("T4").Value = 0.05046448405689858

# Row 5 (Target cluster: MuSCs)
$ws.Range("F5").Value = 0.3333333333333333
$ws.Range("G5").Value = 0.02551366666666667
$ws.Range("H5").Value = 0.076541
$ws.Range("M5").Value = 11.9545335
$ws.Range("N5").Value = 23.909067
$ws.Range("O5").Value = 0.4359947946767024
$ws.Range("P5").Value = 0.3400888529957002
$ws.Range("Q5").Value = 0.3050039828745
$ws.Range("R5").Value = 1.830023897247
$ws.Range("S5").Value = 0.4359947946767024
$ws.Range("T5").Value = 0.3400888529957002

# Row 6 (Target cluster: Resolving-Mac)
$ws.Range("F6").Value = 0.3333333333333333
$ws.Range("G6").Value = 0.02551366666666667
$ws.Range("H6").Value = 0.076541
$ws.Range("M6").Value = 0.5368526666666666
$ws.Range("N6").Value = 1.610558
$ws.Range("O6").Value = 0.01957959866647858
$ws.Range("P6").Value = 0.022909000292778
$ws.Range("Q6").Value = 0.01369707998644444
$ws.Range("R6").Value = 0.123273719878
$ws.Range("S6").Value = 0.01957959866647858
$ws.Range("T6").Value = 0.022909000292778
